$wb = $excel.ActiveWorkbook

# Values for the new "p_eb_c_fix" column (C) and the updated "p_eb_c_inv" column (B)
# per sheet (sheet order matches tab order: 2025, 2030, 2035, 2040, 2045, 2050)
$invValues = @(110000, 105000, 103750, 102500, 101250, 100000)
$fixValues = @(1110, 1080, 1058.0550000000001, 1031.47, 1004.885, 978.3)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Add the new header in C1
    $ws.Range("C1").Value = "p_eb_c_fix"

    # Update the investment cost value in B2
    $ws.Range("B2").Value = $invValues[$i - 1]

    # Set the new fixed cost value in C2
    $ws.Range("C2").Value = $fixValues[$i - 1]
}
